# Update main.py to read and print data from Excel file
# -> regenerate the worksheet data: clear the old name/password header
#    and replace it with a names+numbers table, formatted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old header cells (A1:B1) - the old "name"/"password" labels go away
$ws.Range("A1:B1").ClearContents()

# New data table: name in column A, number in column B, starting row 2
$names   = @("Abdullah","Ahmed","Ali","Hassan","Umar","Bilal","Zain","Saad")
$numbers = @(12345788,84930211,77219834,66549012,90871234,33458790,55671209,78123456)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $numbers[$i]
}

# --- Formatting ---
# Row 1 (now empty) keeps a bold / centered / wrapped look.
# Build the style fully on A1 first (properties are coalesced into one style
# while we keep touching the same cell object), then propagate it with a
# format-only paste so every destination cell shares a single style record.
$a1 = $ws.Cells.Item(1, 1)
$a1.VerticalAlignment = -4108   # xlCenter
$a1.WrapText = $true
$a1.HorizontalAlignment = -4108 # xlCenter
$a1.Font.Bold = $true
$a1.Copy() | Out-Null
$ws.Range("A1:B1").PasteSpecial(-4122) | Out-Null # xlPasteFormats

# Data rows get vertical-centered, wrapped text - same trick, built on A2.
$a2 = $ws.Cells.Item(2, 1)
$a2.VerticalAlignment = -4108   # xlCenter
$a2.WrapText = $true
$a2.Copy() | Out-Null
$ws.Range("A2:B9").PasteSpecial(-4122) | Out-Null # xlPasteFormats

$excel.CutCopyMode = $false

# Column widths (auto-fit based on content)
$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Columns.Item(2).AutoFit() | Out-Null

# Selection: whole first row selected (as in the edited file)
$ws.Range("A1:XFD1").Select() | Out-Null
